$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Write the new text cells first. "minutes/month" must land in the shared
# string table before "mm" so the new unique-string indices come out in the
# same order as the target workbook (34 = minutes/month, 35 = mm).
$ws.Range("H10").Value = "minutes/month"
$ws.Range("A10").Value = "mm"
$ws.Range("C10").Value = "numpy.random"
$ws.Range("D10").Value = "choice"

# Date-formatted cells: clone the existing date style (used by I2/J2/L2,
# cellXfs index 3, numFmtId 14) via a formats-only paste, then overwrite the
# value so the numeric value lands without disturbing the copied style.
$ws.Range("I2").Copy()
$ws.Range("I10").PasteSpecial(-4122)
$ws.Range("I10").Value = 42005

$ws.Range("J2").Copy()
$ws.Range("J10").PasteSpecial(-4122)
$ws.Range("J10").Value = 42370

$ws.Range("L2").Copy()
$ws.Range("L10").PasteSpecial(-4122)
$ws.Range("L10").Value = 42005

# Plain numeric cell, default style.
$ws.Range("K10").Value = 0.5

# Formula cell with a new scientific-notation number format (numFmtId 11).
$ws.Range("E10").Formula = "=6000000000*60"
$ws.Range("E10").NumberFormat = "0.00E+00"
$ws.Range("F10").NumberFormat = "0.00E+00"

# Move the active selection to the new row, matching the saved view state.
$ws.Activate() | Out-Null
$ws.Range("E10").Select() | Out-Null
